$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 0.4115256666666666
$ws.Range("H2").Value = 1.234577
$ws.Range("I2").Value = 0.2245998342667577
$ws.Range("J2").Value = 0.2245998342667577
$ws.Range("M2").Value = 18.43631966666667
$ws.Range("N2").Value = 55.308959
$ws.Range("O2").Value = 0.6034704469962782
$ws.Range("P2").Value = 0.603470446996278
$ws.Range("Q2").Value = 7.587018741704776
$ws.Range("R2").Value = 68.28316867534299
$ws.Range("S2").Value = 0.1355393623802503
$ws.Range("T2").Value = 0.1355393623802502
$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 0.4115256666666666
$ws.Range("H3").Value = 1.234577
$ws.Range("I3").Value = 0.2245998342667577
$ws.Range("J3").Value = 0.2245998342667577
$ws.Range("O3").Value = 0.1750419652256785
$ws.Range("P3").Value = 0.1750419652256784
$ws.Range("Q3").Value = 2.200682199703889
$ws.Range("R3").Value = 19.806139797335
$ws.Range("S3").Value = 0.03931439637941495
$ws.Range("T3").Value = 0.03931439637941495
$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 0.4115256666666666
$ws.Range("H4").Value = 1.234577
$ws.Range("I4").Value = 0.2245998342667577
$ws.Range("J4").Value = 0.2245998342667577
$ws.Range("M4").Value = 6.766555
$ws.Range("N4").Value = 20.299665
$ws.Range("O4").Value = 0.2214875877780434
$ws.Range("P4").Value = 0.2214875877780434
$ws.Range("Q4").Value = 2.784611057411666
$ws.Range("R4").Value = 25.061499516705
$ws.Range("S4").Value = 0.0497460755070925
$ws.Range("T4").Value = 0.04974607550709249
$ws.Range("I5").Value = 0.3944722233087159
$ws.Range("J5").Value = 0.3944722233087159
$ws.Range("M5").Value = 18.43631966666667
$ws.Range("N5").Value = 55.308959
$ws.Range("O5").Value = 0.6034704469962782
$ws.Range("P5").Value = 0.603470446996278
$ws.Range("Q5").Value = 13.32533552883456
$ws.Range("R5").Value = 119.928019759511
$ws.Range("S5").Value = 0.2380523289277264
$ws.Range("T5").Value = 0.2380523289277264
$ws.Range("I6").Value = 0.3944722233087159
$ws.Range("J6").Value = 0.3944722233087159
$ws.Range("O6").Value = 0.1750419652256785
$ws.Range("P6").Value = 0.1750419652256784
$ws.Range("S6").Value = 0.06904919319490033
$ws.Range("T6").Value = 0.06904919319490031
$ws.Range("I7").Value = 0.3944722233087159
$ws.Range("J7").Value = 0.3944722233087159
$ws.Range("M7").Value = 6.766555
$ws.Range("N7").Value = 20.299665
$ws.Range("O7").Value = 0.2214875877780434
$ws.Range("P7").Value = 0.2214875877780434
$ws.Range("Q7").Value = 4.890705812198334
$ws.Range("R7").Value = 44.016352309785
$ws.Range("S7").Value = 0.08737070118608915
$ws.Range("T7").Value = 0.08737070118608914
$ws.Range("G8").Value = 0.6979596666666668
$ws.Range("H8").Value = 2.093879
$ws.Range("I8").Value = 0.3809279424245264
$ws.Range("J8").Value = 0.3809279424245264
$ws.Range("M8").Value = 18.43631966666667
$ws.Range("N8").Value = 55.308959
$ws.Range("O8").Value = 0.6034704469962782
$ws.Range("P8").Value = 0.603470446996278
$ws.Range("Q8").Value = 12.86780752910678
$ws.Range("R8").Value = 115.810267761961
$ws.Range("S8").Value = 0.2298787556883014
$ws.Range("T8").Value = 0.2298787556883014
$ws.Range("G9").Value = 0.6979596666666668
$ws.Range("H9").Value = 2.093879
$ws.Range("I9").Value = 0.3809279424245264
$ws.Range("J9").Value = 0.3809279424245264
$ws.Range("O9").Value = 0.1750419652256785
$ws.Range("P9").Value = 0.1750419652256784
$ws.Range("Q9").Value = 3.732421909393889
$ws.Range("R9").Value = 33.591797184545
$ws.Range("S9").Value = 0.0666783756513632
$ws.Range("T9").Value = 0.06667837565136318
$ws.Range("G10").Value = 0.6979596666666668
$ws.Range("H10").Value = 2.093879
$ws.Range("I10").Value = 0.3809279424245264
$ws.Range("J10").Value = 0.3809279424245264
$ws.Range("M10").Value = 6.766555
$ws.Range("N10").Value = 20.299665
$ws.Range("O10").Value = 0.2214875877780434
$ws.Range("P10").Value = 0.2214875877780434
$ws.Range("Q10").Value = 4.722782472281668
$ws.Range("R10").Value = 42.50504225053501
$ws.Range("S10").Value = 0.08437081108486175
$ws.Range("T10").Value = 0.08437081108486173
